$wb = $excel.ActiveWorkbook

# Sheet "1": update adjustable-cell values (E4/F4, E6/F6).
# D10 (=MMULT(E4:F4,TRANSPOSE(E5:F5))) recalculates automatically.
$ws1 = $wb.Worksheets.Item("1")
$ws1.Range("E4").Value = 3
$ws1.Range("F4").Value = 0
$ws1.Range("E6").Value = 5
$ws1.Range("F6").Value = 10

# Move/extend the visible selection on sheet "1" to match the saved view state.
$ws1.Range("E4:F6").Select() | Out-Null

# Sheet "2": update the two adjusted values.
$ws2 = $wb.Worksheets.Item("2")
$ws2.Range("F4").Value = 1.1599999999999999
$ws2.Range("F5").Value = 20
